# "add mysql connect info": three new "Sql*" config rows (SqlPort, SqlUser,
# SqlPwd) plus renaming the existing "Pwd" row to "SqlIP", mirroring the
# existing ServerID/IP/Port/Pwd rows already on the "Property" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three rows below row 5, copying its formatting/layout each time so
# the new rows keep the same per-cell styles as the existing data rows.
$ws.Rows("5").Copy()
$ws.Rows("6").Insert()

$ws.Rows("5").Copy()
$ws.Rows("7").Insert()

$ws.Rows("5").Copy()
$ws.Rows("8").Insert()

# Row 5 ("Pwd") becomes "SqlIP" (string).
$ws.Range("A5").Value = "SqlIP"
$ws.Range("I5").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Row 6: SqlPort (int)
$ws.Range("A6").Value = "SqlPort"
$ws.Range("I6").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "int"

# Row 7: SqlUser (string)
$ws.Range("A7").Value = "SqlUser"
$ws.Range("B7").Value = "string"

# Row 8: SqlPwd (string)
$ws.Range("A8").Value = "SqlPwd"
$ws.Range("B8").Value = "string"

# Match the new selection recorded in the saved workbook.
$null = $ws.Range("C15").Select()
